# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D sometimes look numeric (e.g. "0.998", "6.81",
# "27.262.03" grouped-thousands style) but must stay plain text, exactly as
# authored, so a leading apostrophe forces text entry the way a user would
# in the Excel UI. Percent cells in column E already carry padding spaces
# that keep them safely text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.262.03"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "'1.650.94"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").Value = "'219.04"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "'20.32"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'1.876.73"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'1.646.38"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'4.15"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'0.543"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Value = "'67.94"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "'27.212.93"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "'0.0₃0739"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'222.59"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "'0.998"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'2.45"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'147.56"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'7.43"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'15.89"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "'0.0508"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'3.05"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").Value = "'1.274.99"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").Value = "'0.546"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").Value = "'0.847"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'5.39"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'63.74"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "'1.786.91"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'92.67"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'0.0₆0103"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'7.71"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'0.0977"
$ws.Range("E51").Value = "  +0.76%  "
